$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RA")

# New rows 12-14: additional CMU contacts (David Held, Deepak Pathak, AV Center)
$ws.Range("A12").Value = "David Held"
$ws.Range("C12").Value = "CMU"
$ws.Range("D12").Value = "https://davheld.github.io/"

$ws.Range("A13").Value = "Deepak Pathak"
$ws.Range("C13").Value = "CMU"
$ws.Range("D13").Value = "https://www.cs.cmu.edu/~dpathak/#ResearchGroup"

$ws.Range("A14").Value = "AV Center"
$ws.Range("C14").Value = "CMU"
$ws.Range("D14").Value = "https://labs.ri.cmu.edu/av-center/we-are-hiring/"

# Row 23: lone marker cell in column F
$ws.Range("F23").Value = " "

# Restore the active sheet/selection to A15
$ws.Activate()
[void]$ws.Range("A15").Select()
